# Updated main GSC export data:
# The oldest date row ("2025-11-20") was dropped from the export and every
# subsequent row moved up one place, so the whole breadcrumb history now
# starts at "2025-11-21" and the table shrinks by one row (C87 instead of
# C88). Deleting row 2 of the "Chart" sheet and letting Excel shift the
# remaining rows up reproduces this exactly (dates, counts, shared-string
# table renumbering, and the header indices used by the "Critical issues"
# and "Non-critical issues" sheets all update automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
